$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping_ind")

# Simplify the steel description: remove the "/RME" segment from the two
# affected lines in the multi-line industrial mapping text (cell B2).
$cell = $ws.Range("B2")
$text = $cell.Value2
$text = $text -replace "31% S/LFM\+CDL/RME/H:1", "31% S/LFM+CDL/H:1"
$text = $text -replace "2% S/LFM\+CDM/RME/H:1", "2% S/LFM+CDM/H:1"
$cell.Value = $text

# Wrap the text and grow the row to fit the full (now slightly shorter)
# multi-line description.
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6

# Leave the same selection state captured in the saved workbook.
$ws.Range("E2:E6").Select()

$wb.Save()
